# Append the next day's GSC export row to the legacy HTTPS/Non-HTTPS
# export on the "Chart" sheet, right after the last existing dated row
# (e.g. "2025-11-17" -> "2025-11-18"). Non-HTTPS/HTTPS URL counts default
# to 0, matching every other freshly-appended day in this export.

function Pad2([int]$n) {
    if ($n -lt 10) {
        return "0" + $n
    }
    return "" + $n
}

function Add-OneDay([int]$y, [int]$m, [int]$d) {
    $daysInMonth = @(31,28,31,30,31,30,31,31,30,31,30,31)
    $isLeap = (($y % 4 -eq 0) -and ($y % 100 -ne 0)) -or ($y % 400 -eq 0)
    if ($isLeap) {
        $daysInMonth[1] = 29
    }

    $d = $d + 1
    if ($d -gt $daysInMonth[$m - 1]) {
        $d = 1
        $m = $m + 1
        if ($m -gt 12) {
            $m = 1
            $y = $y + 1
        }
    }
    return @($y, $m, $d)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Locate the last populated row in the date column (A) and the row right
# after it, where the new day's data belongs.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Read the last date ("yyyy-MM-dd" text) and compute the following day.
$lastDateText = $ws.Cells.Item($lastRow, 1).Text
$parts = $lastDateText.Split("-")
$y = [int]$parts[0]
$m = [int]$parts[1]
$d = [int]$parts[2]

$next = Add-OneDay $y $m $d
$nextDateText = "" + $next[0] + "-" + (Pad2 $next[1]) + "-" + (Pad2 $next[2])

# Write the new date as plain text (not as a converted date serial): enter
# it as a string-literal formula, then collapse the formula down to its
# cached text result via copy / paste-values, matching the plain text
# shared-string cells used by the rest of the date column.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.Formula = '="' + $nextDateText + '"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

# New day starts with zero recorded Non-HTTPS / HTTPS URLs, same as every
# other freshly appended row in this export.
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 0
